$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '98.029.50'
$ws.Range("E2").Value = '  -1.26%  '
$ws.Range("D3").Value = '3.422.49'
$ws.Range("E3").Value = '  +3.19%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '255.38'
$ws.Range("E5").Value = '  -0.80%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '656.77'
$ws.Range("E6").Value = '  +4.39%  '
$ws.Range("E7").Value = '  +2.68%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.435'
$ws.Range("E8").Value = '  +5.73%  '
$ws.Range("E9").Value = '  +8.86%  '
$ws.Range("E10").Value = '  +0.01%  '
$ws.Range("D11").Value = '3.420.02'
$ws.Range("E11").Value = '  +3.18%  '
$ws.Range("E12").Value = '  +4.24%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '42.06'
$ws.Range("E13").Value = '  +2.41%  '
$ws.Range("E14").Value = '  +15.20%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000260'
$ws.Range("E15").Value = '  +2.73%  '
$ws.Range("D16").Value = '97.711.17'
$ws.Range("E16").Value = '  -1.26%  '
$ws.Range("D17").Value = '4.064.69'
$ws.Range("E17").Value = '  +3.30%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.60'
$ws.Range("E18").Value = '  +34.05%  '
$ws.Range("D19").Value = '3.422.18'
$ws.Range("E19").Value = '  +3.23%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.58'
$ws.Range("E20").Value = '  +11.46%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.503'
$ws.Range("E21").Value = '  +48.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.47'
$ws.Range("E22").Value = '  -1.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.75'
$ws.Range("E23").Value = '  +13.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '507.03'
$ws.Range("E24").Value = '  +3.94%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000207'
$ws.Range("E25").Value = '  +0.83%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.23'
$ws.Range("E26").Value = '  +7.80%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '99.24'
$ws.Range("E27").Value = '  +11.33%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.79'
$ws.Range("D29").Value = '3.606.54'
$ws.Range("E29").Value = '  +3.59%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.152'
$ws.Range("E30").Value = '  -1.09%  '
$ws.Range("E31").Value = '  +5.89%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.42'
$ws.Range("E32").Value = '  +6.44%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  -0.24%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  +0.06%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.572'
$ws.Range("E35").Value = '  +17.78%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '29.84'
$ws.Range("E36").Value = '  +6.44%  '
$ws.Range("E37").Value = '  +16.37%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.81'
$ws.Range("E38").Value = '  +5.42%  '
$ws.Range("E39").Value = '  +14.54%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.154'
$ws.Range("E40").Value = '  +1.23%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '524.15'
$ws.Range("E41").Value = '  +5.16%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '24.76'
$ws.Range("E42").Value = '  +0.07%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.872'
$ws.Range("E43").Value = '  +10.54%  '
$ws.Range("E44").Value = '  -4.47%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0421'
$ws.Range("E45").Value = '  +23.71%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.53'
$ws.Range("E46").Value = '  +14.51%  '
$ws.Range("E47").Value = '  +2.19%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.28'
$ws.Range("E48").Value = '  +12.92%  '
$ws.Range("E49").Value = '  +0.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.60'
$ws.Range("E50").Value = '  +13.71%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.07'
$ws.Range("E51").Value = '  +5.00%  '
